$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trip added (Thurgau Säntis Classic) increases the "begleitete Reisen"
# headcount on row 8 (Toskana 2017 sheet) from 19 to 21 participants.
$ws.Range("B8").Value = 21

# Reflect the author's last active cell selection in the saved view state.
$ws.Range("F25").Select() | Out-Null
